# Add a new localization entry ("strChkWindowPosition") to the lv-LV
# translation table. The table (Tabla13) is kept sorted alphabetically by
# the "Key" column (column C), and this new key sorts right before the
# existing "strDifferentiationAlgorithms" row, i.e. it belongs at sheet
# row 34 (the table's data starts at row 3, header is row 2).
#
# Inserting a whole sheet row there shifts every following row down by
# one (old row 34 -> new row 35, ..., old row 192 -> new row 193), which
# matches the rest of the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Push rows 34..192 down to 35..193 and leave row 34 empty (but formatted
# like the surrounding rows, since Insert copies formatting from above).
$ws.Rows("34:34").Insert()

# Fill in the new row's data.
$ws.Range("B34").Value = "localization\strings"
$ws.Range("C34").Value = "strChkWindowPosition"
$ws.Range("D34").Value = 'In "settings" form, tab "User interface"'
$ws.Range("E34").Value = "Remember window position and size on startup"

# Grow the table (ListObject) so it covers the newly inserted row.
$lo.Resize($ws.Range("B2:F193"))
